$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: "수집일" (collection date) header + a constant collection-date
# value (2024-11-26, serial 45622) for every data row (2-112).
$ws.Range("F1").Value = "수집일"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("F1").VerticalAlignment = -4160    # xlTop
$ws.Range("F1").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("F1").Borders.Item(10).LineStyle = 1  # xlEdgeRight

for ($r = 2; $r -le 112; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value = 45622
    $cell.NumberFormat = "mm-dd-yy"
}
